$d = $word.ActiveDocument

# The first paragraph in the body is the title; the second paragraph is an
# empty "spacer" paragraph (bold rPr, sz24/szCs24) that precedes the first
# dated entry ("15.04.2024"). The new version-history entry is inserted
# between the title and that spacer, i.e. a new spacer + new dated entry is
# placed first, followed by the original (untouched) spacer paragraph.
$target = $d.Paragraphs.Item(2).Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>21.05.2024</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Inklusions- und </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Lernstandszeugnisse</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>(Ankreuzberichte)</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="18"/>
    </w:numPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Linien zwischen Kompetenzen im Wahlpflichtbereich wurden erg&#228;nzt.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
'@

$target.InsertXML($xml)

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
Write-Output ("P1: [" + $d.Paragraphs.Item(1).Range.Text + "]")
Write-Output ("P2: [" + $d.Paragraphs.Item(2).Range.Text + "]")
Write-Output ("P3: [" + $d.Paragraphs.Item(3).Range.Text + "]")
Write-Output ("P4: [" + $d.Paragraphs.Item(4).Range.Text + "]")
Write-Output ("P5: [" + $d.Paragraphs.Item(5).Range.Text + "]")
Write-Output ("P6: [" + $d.Paragraphs.Item(6).Range.Text + "]")
Write-Output ("P7: [" + $d.Paragraphs.Item(7).Range.Text + "]")
